# Rotate two blocks of columns left by one position (the first column's
# header+value moves to become the last column of the block), shifting
# everything else one position to the left.
#
# Block 1: AG:AM (columns 33:39) on rows 1 (headers) and 2 (values)
# Block 2: AS:BB (columns 45:54) on rows 1 (headers) and 2 (values)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-ColumnsLeft {
    param($StartCol, $EndCol, $Row)

    $count = $EndCol - $StartCol + 1
    $original = @(0,0,0,0,0,0,0,0,0,0,0,0)

    for ($i = 0; $i -lt $count; $i++) {
        $original[$i] = $ws.Cells.Item($Row, $StartCol + $i).Value()
    }

    for ($i = 0; $i -lt $count; $i++) {
        $srcIndex = ($i + 1) % $count
        $ws.Cells.Item($Row, $StartCol + $i).Value = $original[$srcIndex]
    }
}

# Column AG = 33, AM = 39
Rotate-ColumnsLeft 33 39 1
Rotate-ColumnsLeft 33 39 2

# Column AS = 45, BB = 54
Rotate-ColumnsLeft 45 54 1
Rotate-ColumnsLeft 45 54 2
